$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.326.45"
Set-TextValue $ws.Range("E2") "  +0.21%  "
Set-TextValue $ws.Range("D3") "2.273.94"
Set-TextValue $ws.Range("E3") "  -0.50%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.34%  "
Set-TextValue $ws.Range("D5") "308.84"
Set-TextValue $ws.Range("E5") "  -4.00%  "
Set-TextValue $ws.Range("D6") "102.71"
Set-TextValue $ws.Range("E6") "  +0.38%  "
Set-TextValue $ws.Range("E7") "  -0.36%  "
Set-TextValue $ws.Range("E8") "  -0.11%  "
Set-TextValue $ws.Range("D9") "0.596"
Set-TextValue $ws.Range("E9") "  -1.29%  "
Set-TextValue $ws.Range("D10") "38.29"
Set-TextValue $ws.Range("E10") "  -2.69%  "
Set-TextValue $ws.Range("D11") "0.0892"
Set-TextValue $ws.Range("E11") "  -1.02%  "
Set-TextValue $ws.Range("D12") "8.16"
Set-TextValue $ws.Range("E12") "  -1.49%  "
Set-TextValue $ws.Range("E13") "  +1.02%  "
Set-TextValue $ws.Range("D14") "0.965"
Set-TextValue $ws.Range("E14") "  +0.39%  "
Set-TextValue $ws.Range("D15") "14.91"
Set-TextValue $ws.Range("E15") "  -1.09%  "
Set-TextValue $ws.Range("D16") "2.618.42"
Set-TextValue $ws.Range("E16") "  -0.55%  "
Set-TextValue $ws.Range("D17") "2.274.03"
Set-TextValue $ws.Range("E17") "  -0.28%  "
Set-TextValue $ws.Range("D18") "42.275.80"
Set-TextValue $ws.Range("E18") "  -0.11%  "
Set-TextValue $ws.Range("D19") "7.19"
Set-TextValue $ws.Range("E19") "  -2.03%  "
Set-TextValue $ws.Range("E20") "  -1.60%  "
Set-TextValue $ws.Range("D21") "12.89"
Set-TextValue $ws.Range("E21") "  +1.10%  "
Set-TextValue $ws.Range("D22") "72.46"
Set-TextValue $ws.Range("E22") "  -0.66%  "
Set-TextValue $ws.Range("D23") "3.35"
Set-TextValue $ws.Range("E23") "  -6.31%  "
Set-TextValue $ws.Range("D24") "260.86"
Set-TextValue $ws.Range("E24") "  -2.52%  "
Set-TextValue $ws.Range("D25") "2.15"
Set-TextValue $ws.Range("E25") "  -3.24%  "
Set-TextValue $ws.Range("E26") "  +0.49%  "
Set-TextValue $ws.Range("E27") "  -2.52%  "
Set-TextValue $ws.Range("D28") "2.27"
Set-TextValue $ws.Range("E28") "  -1.99%  "
Set-TextValue $ws.Range("D29") "6.79"
Set-TextValue $ws.Range("E29") "  +12.56%  "
Set-TextValue $ws.Range("D30") "22.00"
Set-TextValue $ws.Range("E30") "  -2.17%  "
Set-TextValue $ws.Range("E31") "  -6.95%  "
Set-TextValue $ws.Range("D32") "163.30"
Set-TextValue $ws.Range("E32") "  -0.41%  "
Set-TextValue $ws.Range("D33") "0.0846"
Set-TextValue $ws.Range("E33") "  -2.77%  "
Set-TextValue $ws.Range("E34") "  -2.89%  "
Set-TextValue $ws.Range("E35") "  +1.16%  "
Set-TextValue $ws.Range("E36") "  -3.48%  "
Set-TextValue $ws.Range("D37") "4.45"
Set-TextValue $ws.Range("E37") "  -2.95%  "
Set-TextValue $ws.Range("E38") "  -3.23%  "
Set-TextValue $ws.Range("E39") "  -1.19%  "
Set-TextValue $ws.Range("D40") "2.71"
Set-TextValue $ws.Range("E40") "  -1.42%  "
Set-TextValue $ws.Range("E41") "  +0.75%  "
Set-TextValue $ws.Range("D42") "98.92"
Set-TextValue $ws.Range("E42") "  +9.32%  "
Set-TextValue $ws.Range("E43") "  -0.34%  "
Set-TextValue $ws.Range("D44") "67.89"
Set-TextValue $ws.Range("E44") "  -0.71%  "
Set-TextValue $ws.Range("D45") "0.223"
Set-TextValue $ws.Range("E45") "  -0.82%  "
Set-TextValue $ws.Range("D46") "1.710.23"
Set-TextValue $ws.Range("E46") "  +6.99%  "
Set-TextValue $ws.Range("D47") "11.77"
Set-TextValue $ws.Range("E47") "  -3.11%  "
Set-TextValue $ws.Range("D48") "108.87"
Set-TextValue $ws.Range("E48") "  -4.03%  "
Set-TextValue $ws.Range("D49") "75.25"
Set-TextValue $ws.Range("E49") "  -5.44%  "
Set-TextValue $ws.Range("D50") "8.55"
Set-TextValue $ws.Range("E50") "  -4.37%  "
Set-TextValue $ws.Range("D51") "5.08"
Set-TextValue $ws.Range("E51") "  -2.55%  "
